$wb = $excel.ActiveWorkbook

# Sheet: ALC (index 1)
$ws = $wb.Worksheets.Item(1)
$ws.Range("H40").Value = 41668670
$ws.Range("I40").Value = 1433.3334
$ws.Range("J40").Value = 47621132
$ws.Range("K40").Value = 1433.3334
$ws.Range("L40").Value = 47621132
$ws.Range("M40").Value = -1258.3334
$ws.Range("N40").Value = -47621482
$ws.Range("H43").Value = 2001440.6
$ws.Range("J43").Value = 3334134
$ws.Range("L43").Value = 3334134
$ws.Range("N43").Value = -3334272
$ws.Range("H116").Value = 3295.1177
$ws.Range("I116").Value = 3542.0833
$ws.Range("K116").Value = 3542.0833
$ws.Range("M116").Value = -100.0832999999998
$ws.Range("H134").Value = 67161.336
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 67161.336
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 67161.336
$ws.Range("M134").ClearContents()
$ws.Range("N134").Value = -77301.336
$ws.Range("H137").Value = 9092428
$ws.Range("I137").Value = 1473.7222
$ws.Range("J137").Value = 26317396
$ws.Range("K137").Value = 4421.1666
$ws.Range("L137").Value = 78952188
$ws.Range("M137").Value = -1871.1666
$ws.Range("N137").Value = -78957288
$ws.Range("H138").Value = 4136.242
$ws.Range("I138").Value = 2348.75
$ws.Range("J138").Value = 4293.385
$ws.Range("K138").Value = 7046.25
$ws.Range("L138").Value = 12880.155
$ws.Range("M138").Value = -1906.25
$ws.Range("N138").Value = -23160.155

# Sheet: ARM (index 2)
$ws = $wb.Worksheets.Item(2)
$ws.Range("H32").Value = 19832.037
$ws.Range("I32").Value = 17152.898
$ws.Range("K32").Value = 17152.898
$ws.Range("M32").Value = -16865.898
$ws.Range("H45").Value = 1845.8846
$ws.Range("I45").Value = 1025.9474
$ws.Range("J45").Value = 4071.4285
$ws.Range("K45").Value = 1025.9474
$ws.Range("L45").Value = 4071.4285
$ws.Range("M45").Value = -648.9474
$ws.Range("N45").Value = -4825.4285
$ws.Range("H61").Value = 1504.973
$ws.Range("I61").Value = 1308.1666
$ws.Range("J61").Value = 1868.3077
$ws.Range("K61").Value = 1308.1666
$ws.Range("L61").Value = 1868.3077
$ws.Range("M61").Value = -1096.1666
$ws.Range("N61").Value = -2292.3077
$ws.Range("H74").Value = 937.08887
$ws.Range("I74").Value = 1194.0952
$ws.Range("J74").Value = 712.2083
$ws.Range("K74").Value = 1194.0952
$ws.Range("L74").Value = 712.2083
$ws.Range("M74").Value = -320.0952
$ws.Range("N74").Value = -2460.2083
$ws.Range("H77").Value = 937.08887
$ws.Range("I77").Value = 1194.0952
$ws.Range("J77").Value = 712.2083
$ws.Range("K77").Value = 5970.476
$ws.Range("L77").Value = 3561.0415
$ws.Range("M77").Value = -1602.476
$ws.Range("N77").Value = -12297.0415
$ws.Range("H136").Value = 1504.973
$ws.Range("I136").Value = 1308.1666
$ws.Range("J136").Value = 1868.3077
$ws.Range("K136").Value = 3924.4998
$ws.Range("L136").Value = 5604.9231
$ws.Range("M136").Value = -1374.4998
$ws.Range("N136").Value = -10704.9231

# Sheet: BSM (index 3)
$ws = $wb.Worksheets.Item(3)
$ws.Range("H99").Value = 62502176
$ws.Range("I99").Value = 100002080
$ws.Range("J99").Value = 2337
$ws.Range("K99").Value = 100002080
$ws.Range("L99").Value = 2337
$ws.Range("M99").Value = -100000582
$ws.Range("N99").Value = -5333
$ws.Range("H113").Value = 24659.8
$ws.Range("I113").Value = 24659.8
$ws.Range("K113").Value = 24659.8
$ws.Range("M113").Value = -22489.8

# Sheet: CRP (index 4)
$ws = $wb.Worksheets.Item(4)
$ws.Range("H7").Value = 169.76471
$ws.Range("I7").Value = 153.38461
$ws.Range("K7").Value = 153.38461
$ws.Range("M7").Value = -40.38461000000001
$ws.Range("H22").Value = 283.1111
$ws.Range("I22").Value = 249.71428
$ws.Range("K22").Value = 249.71428
$ws.Range("M22").Value = 100.28572
$ws.Range("H31").Value = 2043.6774
$ws.Range("I31").Value = 956.875
$ws.Range("J31").Value = 2421.6956
$ws.Range("K31").Value = 956.875
$ws.Range("L31").Value = 2421.6956
$ws.Range("M31").Value = -661.875
$ws.Range("N31").Value = -3011.6956
$ws.Range("H34").Value = 2043.6774
$ws.Range("I34").Value = 956.875
$ws.Range("J34").Value = 2421.6956
$ws.Range("K34").Value = 956.875
$ws.Range("L34").Value = 2421.6956
$ws.Range("M34").Value = -754.875
$ws.Range("N34").Value = -2825.6956
$ws.Range("H99").Value = 2167.6
$ws.Range("I99").Value = 1914.5
$ws.Range("J99").Value = 3180
$ws.Range("K99").Value = 1914.5
$ws.Range("L99").Value = 3180
$ws.Range("M99").Value = -416.5
$ws.Range("N99").Value = -6176
$ws.Range("H126").Value = 2167.6
$ws.Range("I126").Value = 1914.5
$ws.Range("J126").Value = 3180
$ws.Range("K126").Value = 5743.5
$ws.Range("L126").Value = 9540
$ws.Range("M126").Value = -3273.5
$ws.Range("N126").Value = -14480

# Sheet: CUL (index 5)
$ws = $wb.Worksheets.Item(5)
$ws.Range("H2").Value = 2941317.8
$ws.Range("I2").Value = 5882501
$ws.Range("J2").Value = 134.4
$ws.Range("K2").Value = 35295006
$ws.Range("L2").Value = 806.4000000000001
$ws.Range("M2").Value = -35294893
$ws.Range("N2").Value = -1032.4
$ws.Range("H3").Value = 3860.3809
$ws.Range("I3").Value = 1897.8572
$ws.Range("J3").Value = 7785.4287
$ws.Range("K3").Value = 5693.571599999999
$ws.Range("L3").Value = 23356.2861
$ws.Range("M3").Value = -5581.571599999999
$ws.Range("N3").Value = -23580.2861
$ws.Range("H20").Value = 3981.5
$ws.Range("I20").Value = 999
$ws.Range("J20").Value = 4578
$ws.Range("K20").Value = 2997
$ws.Range("L20").Value = 13734
$ws.Range("M20").Value = -2770
$ws.Range("N20").Value = -14188
$ws.Range("H113").Value = 175997.38
$ws.Range("I113").Value = 400
$ws.Range("J113").Value = 179133.03
$ws.Range("K113").Value = 1200
$ws.Range("L113").Value = 537399.09
$ws.Range("M113").Value = 970
$ws.Range("N113").Value = -541739.09
$ws.Range("H121").Value = 828.3333
$ws.Range("I121").Value = 330
$ws.Range("J121").Value = 1326.6666
$ws.Range("K121").Value = 990
$ws.Range("L121").Value = 3979.9998
$ws.Range("M121").Value = 320
$ws.Range("N121").Value = -6599.9998
$ws.Range("H122").Value = 370.42856
$ws.Range("I122").Value = 326.51852
$ws.Range("J122").Value = 518.625
$ws.Range("K122").Value = 2938.66668
$ws.Range("L122").Value = 4667.625
$ws.Range("M122").Value = -488.6666800000003
$ws.Range("N122").Value = -9567.625
$ws.Range("H131").Value = 33946.594
$ws.Range("I131").Value = 112755.555
$ws.Range("J131").Value = 3108.3044
$ws.Range("K131").Value = 338266.665
$ws.Range("L131").Value = 9324.913199999999
$ws.Range("M131").Value = -333226.665
$ws.Range("N131").Value = -19404.9132

# Sheet: GSM (index 6)
$ws = $wb.Worksheets.Item(6)
$ws.Range("H113").Value = 1222.7333
$ws.Range("I113").Value = 875.2222
$ws.Range("J113").Value = 1744
$ws.Range("K113").Value = 875.2222
$ws.Range("L113").Value = 1744
$ws.Range("M113").Value = 1294.7778
$ws.Range("N113").Value = -6084
$ws.Range("H122").Value = 482461.8
$ws.Range("I122").Value = 591616.4399999999
$ws.Range("J122").Value = 2181.6
$ws.Range("K122").Value = 1774849.32
$ws.Range("L122").Value = 6544.799999999999
$ws.Range("M122").Value = -1772399.32
$ws.Range("N122").Value = -11444.8

# Sheet: LTW (index 7)
$ws = $wb.Worksheets.Item(7)
$ws.Range("H7").Value = 1455.6
$ws.Range("I7").Value = 1194.5
$ws.Range("J7").Value = 2500
$ws.Range("K7").Value = 1194.5
$ws.Range("L7").Value = 2500
$ws.Range("M7").Value = -1082.5
$ws.Range("N7").Value = -2724
$ws.Range("H40").Value = 249292.14
$ws.Range("I40").Value = 282699.6
$ws.Range("J40").Value = 2077
$ws.Range("K40").Value = 282699.6
$ws.Range("L40").Value = 2077
$ws.Range("M40").Value = -282563.6
$ws.Range("N40").Value = -2349
$ws.Range("H55").Value = 272.6111
$ws.Range("I55").Value = 130.63637
$ws.Range("J55").Value = 495.7143
$ws.Range("K55").Value = 130.63637
$ws.Range("L55").Value = 495.7143
$ws.Range("M55").Value = 42.36363
$ws.Range("N55").Value = -841.7143
$ws.Range("H126").Value = 1455.6
$ws.Range("I126").Value = 1194.5
$ws.Range("J126").Value = 2500
$ws.Range("K126").Value = 3583.5
$ws.Range("L126").Value = 7500
$ws.Range("M126").Value = -1113.5
$ws.Range("N126").Value = -12440

# Sheet: WVR (index 8)
$ws = $wb.Worksheets.Item(8)
$ws.Range("H96").Value = 2321.3684
$ws.Range("I96").Value = 2277.5386
$ws.Range("K96").Value = 2277.5386
$ws.Range("M96").Value = -904.5385999999999
$ws.Range("H132").Value = 1398.3334
$ws.Range("I132").Value = 1081.6383
$ws.Range("J132").Value = 2543.3076
$ws.Range("K132").Value = 3244.9149
$ws.Range("L132").Value = 7629.9228
$ws.Range("M132").Value = -714.9149000000002
$ws.Range("N132").Value = -12689.9228
$ws.Range("H136").Value = 2568.56
$ws.Range("I136").Value = 3411.1667
$ws.Range("J136").Value = 1790.7693
$ws.Range("K136").Value = 10233.5001
$ws.Range("L136").Value = 5372.3079
$ws.Range("M136").Value = -7683.500100000001
$ws.Range("N136").Value = -10472.3079
